# Apply a 3-way rotation of record data among rows 2, 3 and 5:
#   new row2 <- old row5
#   new row3 <- old row2
#   new row5 <- old row3
# Row 4 and all other rows/columns are left untouched.
#
# Target values are hard-coded (rather than read back from the sheet and
# re-written) so that numeric values are written with the exact same
# textual precision intended by the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (becomes what used to be row 5) ----
$ws.Range("A2").Value2 = 111801751
$ws.Range("B2").Value2 = 89405
$ws.Range("D2").Value2 = "NT"
$ws.Range("E2").Value2 = 1202
$ws.Range("F2").Value2 = "Ullticka"
$ws.Range("G2").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H2").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q2").Value2 = 507350.4792540679
$ws.Range("R2").Value2 = 6946859.423841522

# ---- Row 3 (becomes what used to be row 2) ----
$ws.Range("A3").Value2 = 111801756
$ws.Range("B3").Value2 = 77267
$ws.Range("D3").Value2 = "NT"
$ws.Range("E3").Value2 = 6446
$ws.Range("F3").Value2 = "Kolflarnlav"
$ws.Range("G3").Value2 = "Carbonicola anthracophila"
$ws.Range("H3").Value2 = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q3").Value2 = 507338.8401105705
$ws.Range("R3").Value2 = 6946917.328540379

# ---- Row 5 (becomes what used to be row 3) ----
$ws.Range("A5").Value2 = 111801760
$ws.Range("B5").Value2 = 89558
$ws.Range("D5").Value2 = "VU"
$ws.Range("E5").Value2 = 1503
$ws.Range("F5").Value2 = "Gräddporing"
$ws.Range("G5").Value2 = "Sidera lenis"
$ws.Range("H5").Value2 = "(P.Karst.) Miettinen"
$ws.Range("Q5").Value2 = 507292.6252952328
$ws.Range("R5").Value2 = 6946995.844692842
